# Add "Flow vs R1L" summary block (rows 38-40, cols F:J) to Sheet1,
# mirroring the existing Kpl summary block at rows 22-24 (cols B:F) but
# built from column F ("Flow_Lac") instead of column B ("Kpl").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (r=38): same condition labels as the existing summary header
$ws.Range("G38").Value = "HK-2"
$ws.Range("H38").Value = "UMRC6"
$ws.Range("I38").Value = "UOK262"
$ws.Range("J38").Value = "UOK + DIDS"

# Row 39: row label + AVERAGE of the F column per condition block
$ws.Range("F39").Value = "Flow_Lac"
$ws.Range("G39").Formula = "=AVERAGE(F`$1:F`$3)"
$ws.Range("H39").Formula = "=AVERAGE(F`$4:F`$6)"
$ws.Range("I39").Formula = "=AVERAGE(F`$9:F`$11)"
$ws.Range("J39").Formula = "=AVERAGE(F`$13:F`$16)"

# Row 40: standard error of the mean for the same F column blocks
$ws.Range("G40").Formula = "=STDEV(F`$1:F`$3)/SQRT(COUNT(F`$1:F`$3))"
$ws.Range("H40").Formula = "=STDEV(F`$4:F`$6)/SQRT(COUNT(F`$4:F`$6))"
$ws.Range("I40").Formula = "=STDEV(F`$9:F`$11)/SQRT(COUNT(F`$9:F`$11))"
$ws.Range("J40").Formula = "=STDEV(F`$13:F`$16)/SQRT(COUNT(F`$13:F`$16))"

# Match the author's final selection on the newly added block
$ws.Range("F38:J40").Select()
